$d = $word.ActiveDocument

$replacements = @(
    @("346×8=2768", "641×8=5128"),
    @("225×8=1800", "657×9=5913"),
    @("874×3=2622", "772×2=1544"),
    @("612×9=5508", "997×2=1994"),
    @("640×9=5760", "353×2=706"),
    @("564×6=3384", "129×8=1032"),
    @("201×5=1005", "114×3=342"),
    @("366×2=732", "507×3=1521"),
    @("548×8=4384", "538×4=2152"),
    @("866×5=4330", "291×4=1164"),
    @("527×2=1054", "652×7=4564"),
    @("435×9=3915", "619×5=3095"),
    @("261×3=783", "429×5=2145"),
    @("607×8=4856", "580×5=2900"),
    @("876×4=3504", "950×5=4750"),
    @("112×7=784", "393×9=3537"),
    @("476×7=3332", "932×5=4660"),
    @("333×5=1665", "930×8=7440"),
    @("838×7=5866", "463×3=1389"),
    @("132×2=264", "196×2=392"),
    @("868×6=5208", "136×5=680"),
    @("430×3=1290", "534×9=4806"),
    @("590×8=4720", "643×8=5144"),
    @("610×9=5490", "457×6=2742"),
    @("908×3=2724", "326×2=652")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
